$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 4954
$ws1.Range("F5").Value = 32
$ws1.Range("F6").Value = 32
$ws1.Range("G6").Value = 55
$ws1.Range("F8").Value = 493

# Sheet "全部类型" (all types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 4954
$ws4.Range("F6").Value = 32
$ws4.Range("F7").Value = 32
$ws4.Range("G7").Value = 55
$ws4.Range("F10").Value = 493
